# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 12 (pushing the
# existing rows 12-31 down to 13-32), and the new row is populated with
# the latest reading for "Ajo" (Agrícola del Norte S.A. de Arica).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 12, shifting rows 12-31 down to 13-32.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with this week's reading.
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 44925
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112003
$ws.Cells.Item(12, 7).Value = "Ajo"
$ws.Cells.Item(12, 8).Value = "Chino"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 250
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 14600
$ws.Cells.Item(12, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(12, 15).Value = "China"
$ws.Cells.Item(12, 16).Value = 1460
$ws.Cells.Item(12, 17).Value = 10
$ws.Cells.Item(12, 18).Value = "Hortaliza"
